{"js": "// Technical Skills section updates:\n//  - remove \", Python\" from the Languages line\n//  - remove \", Anjular.js\" from the Web Designing line\n//  - add \"Hibernate, \" before \"Log4j,\" on the Standards and Frameworks line\nconst body = context.document.body;\n\n// 1. Remove \", Python\" (Languages: Java, R, SQL, PL/SQL, Python)\nconst pythonResults = body.search(\", Python\", { matchCase: true });\npythonResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < pythonResults.items.length; i++) {\n  pythonResults.items[i].delete();\n}\n\n// 2. Remove \", Anjular.js\" (Web Designing: ... AJAX, Anjular.js)\nconst angularResults = body.search(\", Anjular.js\", { matchCase: true });\nangularResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < angularResults.items.length; i++) {\n  angularResults.items[i].delete();\n}\n\n// 3. Insert \"Hibernate, \" right before \"Log4j,\" (Standards and Frameworks line)\nconst log4jResults = body.search(\"Log4j,\", { matchCase: true });\nlog4jResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < log4jResults.items.length; i++) {\n  log4jResults.items[i].insertText(\"Hibernate, \", Word.InsertLocation.before);\n}\n\nawait context.sync();\n", "ps1": "# Technical Skills section updates:\n#  - remove \", Python\" from the Languages line\n#  - remove \", Anjular.js\" from the Web Designing line\n#  - add \"Hibernate, \" before \"Log4j,\" on the Standards and Frameworks line\n$d = $word.ActiveDocument\n\nfunction Remove-AllOccurrences($doc, $searchText) {\n  $rng = $doc.Content\n  $find = $rng.Find\n  $find.Text = $searchText\n  $find.MatchCase = $true\n  $find.Forward = $true\n  while ($find.Execute()) {\n    $rng.Delete()\n    $rng.Start = $rng.End\n    $rng.End = $doc.Content.End\n  }\n}\n\n# 1. Remove \", Python\" (Languages: Java, R, SQL, PL/SQL, Python)\nRemove-AllOccurrences $d \", Python\"\n\n# 2. Remove \", Anjular.js\" (Web Designing: ... AJAX, Anjular.js)\nRemove-AllOccurrences $d \", Anjular.js\"\n\n# 3. Insert \"Hibernate, \" right before \"Log4j,\" (Standards and Frameworks line)\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.Text = \"Log4j,\"\n$find3.MatchCase = $true\n$find3.Forward = $true\nwhile ($find3.Execute()) {\n  $rng3.InsertBefore(\"Hibernate, \")\n  $rng3.Collapse(0)\n  $rng3.End = $d.Content.End\n}\n"}
